$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 13:05"

# Zambia (row 116): Casos totales, Nuevos casos, Recuperados
$ws.Range("B116").Value = 772
$ws.Range("C116").Value = 11
$ws.Range("E116").Value = 573

# San Marino (row 122)
$ws.Range("B122").Value = 655
$ws.Range("C122").Value = 1
$ws.Range("D122").Value = 211
$ws.Range("E122").Value = 403

# Malta (row 125)
$ws.Range("B125").Value = 569
$ws.Range("C125").Value = 11
$ws.Range("D125").Value = 460
$ws.Range("E125").Value = 103

# Nepal overtakes Estado de Palestina in total cases, so the two countries
# swap rows (134/135) while keeping the list sorted by Casos totales.
$ws.Range("A134").Value = "Nepal"
$ws.Range("A135").Value = "Estado de Palestina"

$ws.Range("B134").Value = 402
$ws.Range("C134").Value = 27
$ws.Range("D134").Value = 37
$ws.Range("E134").Value = 363

$ws.Range("B135").Value = 388
$ws.Range("C135").Value = 0
$ws.Range("D135").Value = 337
$ws.Range("E135").Value = 49

# Uganda (row 148)
$ws.Range("B148").Value = 260
$ws.Range("C148").Value = 12
$ws.Range("D148").Value = 197
$ws.Range("E148").Value = 197
